$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the three changed data values in row 4 (E4, G4, H4)
$ws.Range("E4").Value = 5
$ws.Range("G4").Value = -3
$ws.Range("H4").Value = 13

# Update the active cell / selection to E4
$ws.Range("E4").Select()
